# Updates the cryptos list (Price / Volume(1h) columns) to reflect refreshed
# market data, matching a GitHub Actions scheduled data-refresh commit.
# Numeric-looking Price values are prefixed with a leading apostrophe so
# Excel keeps them as literal text (e.g. "554.85") instead of auto-converting
# them to numbers / floats.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.961.67"
$ws.Range("E2").Value = "  +0.92%  "

$ws.Range("D3").Value = "3.365.32"
$ws.Range("E3").Value = "  +0.75%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'554.85"
$ws.Range("E5").Value = "  +0.64%  "

$ws.Range("D6").Value = "'173.92"
$ws.Range("E6").Value = "  -0.34%  "

$ws.Range("E7").Value = "  +2.52%  "

$ws.Range("D8").Value = "3.354.26"
$ws.Range("E8").Value = "  +0.55%  "

$ws.Range("E9").Value = "  +0.00%  "

$ws.Range("D10").Value = "'0.173"
$ws.Range("E10").Value = "  +6.13%  "

$ws.Range("D11").Value = "'0.637"
$ws.Range("E11").Value = "  +1.70%  "

$ws.Range("D12").Value = "'53.56"
$ws.Range("E12").Value = "  -1.44%  "

$ws.Range("E13").Value = "  +3.61%  "

$ws.Range("E14").Value = "  +1.16%  "

$ws.Range("D15").Value = "3.906.06"
$ws.Range("E15").Value = "  +0.84%  "

$ws.Range("E16").Value = "  +2.08%  "

$ws.Range("D17").Value = "'18.22"
$ws.Range("E17").Value = "  -0.17%  "

$ws.Range("D18").Value = "3.373.40"
$ws.Range("E18").Value = "  +1.06%  "

$ws.Range("D19").Value = "65.074.60"
$ws.Range("E19").Value = "  +1.22%  "

$ws.Range("D20").Value = "'11.88"
$ws.Range("E20").Value = "  +1.57%  "

$ws.Range("E21").Value = "  +1.90%  "

$ws.Range("D22").Value = "'454.86"
$ws.Range("E22").Value = "  +3.94%  "

$ws.Range("D23").Value = "'4.92"
$ws.Range("E23").Value = "  -0.96%  "

$ws.Range("D24").Value = "'4.05"
$ws.Range("E24").Value = "  +0.05%  "

$ws.Range("D25").Value = "'87.30"
$ws.Range("E25").Value = "  +3.49%  "

$ws.Range("D26").Value = "'13.65"
$ws.Range("E26").Value = "  +2.33%  "

$ws.Range("D27").Value = "'10.71"
$ws.Range("E27").Value = "  -0.25%  "

$ws.Range("E28").Value = "  +1.67%  "

$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").Value = "'8.65"
$ws.Range("E29").Value = "  -0.62%  "

$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'31.13"
$ws.Range("E30").Value = "  +4.75%  "

$ws.Range("D31").Value = "'6.53"
$ws.Range("E31").Value = "  -1.62%  "

$ws.Range("D32").Value = "'62.88"
$ws.Range("E32").Value = "  +7.85%  "

$ws.Range("E33").Value = "  -0.04%  "

$ws.Range("D34").Value = "'575.34"
$ws.Range("E34").Value = "  -0.21%  "

$ws.Range("E35").Value = "  +0.16%  "

$ws.Range("E36").Value = "  +0.04%  "

$ws.Range("D37").Value = "'3.66"
$ws.Range("E37").Value = "  +4.18%  "

$ws.Range("E38").Value = "  +0.32%  "

$ws.Range("D39").Value = "'35.63"
$ws.Range("E39").Value = "  +0.54%  "

$ws.Range("D40").Value = "'0.370"
$ws.Range("E40").Value = "  +1.39%  "

$ws.Range("E41").Value = "  -0.91%  "

$ws.Range("D42").Value = "3.071.31"
$ws.Range("E42").Value = "  -0.81%  "

$ws.Range("D43").Value = "'0.0415"
$ws.Range("E43").Value = "  +1.93%  "

$ws.Range("D44").Value = "'2.75"
$ws.Range("E44").Value = "  -0.81%  "

$ws.Range("E45").Value = "  +2.85%  "

$ws.Range("E46").Value = "  -0.38%  "

$ws.Range("E47").Value = "  -1.96%  "

$ws.Range("E48").Value = "  +0.06%  "

$ws.Range("D49").Value = "'141.90"
$ws.Range("E49").Value = "  +4.06%  "

$ws.Range("E50").Value = "  -1.94%  "

$ws.Range("D51").Value = "'8.26"
$ws.Range("E51").Value = "  -0.20%  "
